$d = $word.ActiveDocument

# 1) "cadastrados os quais" -> "cadastrados aos quais"
$oldPara1 = "realizadas por revendedores cadastrados os quais, através de catálogos (revistas), vendem de " + [char]8220 + "porta em porta" + [char]8221
$newPara1 = "realizadas por revendedores cadastrados aos quais, através de catálogos (revistas), vendem de " + [char]8220 + "porta em porta" + [char]8221
$d.Content.Find.Execute(
    $oldPara1,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newPara1,
    2
)

# 2) add comma after "microempresa"
$d.Content.Find.Execute(
    "sua microempresa já que",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sua microempresa, já que",
    2
)

# 3) "(Ponto de Venda)" -> "(ponto de venda)"
$d.Content.Find.Execute(
    "(Ponto de Venda)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(ponto de venda)",
    2
)

# 4) "redução de gastos." -> "redução de custos."
$d.Content.Find.Execute(
    "redução de gastos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "redução de custos.",
    2
)

# 5) " sem problemas." -> "."
$d.Content.Find.Execute(
    "entregar um software para o cliente sem problemas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "entregar um software para o cliente.",
    2
)

Write-Host "Text replacements done"
